$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns carrying new values for each data row (G/H untouched; A is bus index).
$cols = @("B","C","D","E","F","I","J","K","L","M","N")

# One row of new vm_pu values per data row (rows 2..25), aligned to $cols above.
$data = @(
    @(1.02, 1.051327616492103, 1.057001391769153, 1.055005077837572, 1.064842128277711, 1.045860508727283, 1.056356072460677, 1.059737558029918, 1.057746727003507, 1.067556975148261, 1.022491572486312),
    @(1.02, 1.052494811593085, 1.057955194014173, 1.056122298255516, 1.065985567628592, 1.046209669369522, 1.057172245307255, 1.060504754297706, 1.058676527299321, 1.068514876785215, 1.022772460480003),
    @(1.02, 1.053249939910301, 1.058572238221572, 1.056845426037172, 1.066725727466742, 1.046434332139521, 1.057699689611451, 1.061000429027795, 1.059277793150013, 1.069134395142485, 1.022953768284958),
    @(1.02, 1.053567367335855, 1.058831613064274, 1.057149480417672, 1.067036958040134, 1.046528477224533, 1.057921266515857, 1.061208631091416, 1.059530475788334, 1.069394767620412, 1.02302988352844),
    @(1.02, 1.053620663165659, 1.05887516146492, 1.057200535528099, 1.067089218997086, 1.046544266830882, 1.057958460851068, 1.06124357861922, 1.059572897087383, 1.069438481042732, 1.023042657366993),
    @(1.02, 1.053254181501398, 1.058575704120307, 1.056849488623425, 1.066729885879146, 1.04643559130152, 1.057702650963379, 1.061003211738919, 1.059281169858615, 1.069137874538381, 1.022954785759141),
    @(1.02, 1.05172210225568, 1.057323760838637, 1.055382604426728, 1.065228501852431, 1.045978771829607, 1.056632041894915, 1.059996991765978, 1.058061035906462, 1.067880766703347, 1.022586592203989),
    @(1.02, 1.049021363063224, 1.055116661615341, 1.052799348793757, 1.062584965402987, 1.04516407621411, 1.054740297391627, 1.058218114019024, 1.055908081499276, 1.065663190594871, 1.021934371430367),
    @(1.02, 1.047220094499922, 1.053644545459211, 1.051078186251793, 1.060823962520831, 1.044614389257994, 1.053475595306247, 1.057028264443794, 1.054470763342813, 1.064183151135706, 1.02149725287779),
    @(1.02, 1.04643992052577, 1.053006924475879, 1.050333129704329, 1.060061737453761, 1.04437480863045, 1.052927114843271, 1.056512104762732, 1.053847898769691, 1.063541873707439, 1.021307427064008),
    @(1.02, 1.046150095313479, 1.052770055032101, 1.050056414069533, 1.059778657119713, 1.044285582468237, 1.052723255109011, 1.056320236911199, 1.053616463438976, 1.063303611953976, 1.021236834367634),
    @(1.02, 1.046212265368002, 1.052820865640367, 1.050115769105193, 1.059839376840752, 1.044304732436408, 1.052766989598675, 1.0563613997095, 1.053666110546538, 1.063354722766436, 1.021251980497208),
    @(1.02, 1.046415964185, 1.052987345362425, 1.050310255667395, 1.060038337054513, 1.044367437969931, 1.052910266385319, 1.056496247831944, 1.053828769803744, 1.06352218020318, 1.021301593542852),
    @(1.02, 1.046541465205591, 1.053089915196337, 1.050430089399541, 1.060160928802483, 1.044406041731008, 1.052998526728512, 1.05657931321337, 1.053928979476434, 1.063625347907639, 1.021332150781244),
    @(1.02, 1.047271867432374, 1.053686858295931, 1.05112763769229, 1.060874555123528, 1.044630256471883, 1.053511978021454, 1.057062500234238, 1.054512090259921, 1.064225701864171, 1.021509839373016),
    @(1.02, 1.047729971146171, 1.054061254823704, 1.051565248848916, 1.061322273749387, 1.04477048176757, 1.053833822507396, 1.057365336527932, 1.054877726572901, 1.064602177356137, 1.021621151169803),
    @(1.02, 1.047997154886728, 1.054279616192707, 1.051820521168404, 1.061583449552573, 1.044852122080049, 1.054021466418634, 1.057541884506273, 1.055090948235637, 1.06482172959906, 1.02168602434769),
    @(1.02, 1.048088254234012, 1.05435406878406, 1.051907566055172, 1.06167250872494, 1.044879933765443, 1.054085434176159, 1.057602067344179, 1.055163643236923, 1.064896584565937, 1.021708135423068),
    @(1.02, 1.047680823077136, 1.054021087465592, 1.051518295135912, 1.06127423477613, 1.044755452514622, 1.053799300184047, 1.057332854510356, 1.054838502210297, 1.064561789196407, 1.021609213962365),
    @(1.02, 1.04635598091037, 1.052938322030424, 1.050252983361573, 1.059979747025698, 1.04434897925122, 1.052868078530659, 1.056456542380095, 1.053780872805296, 1.063472869903089, 1.021286986028697),
    @(1.02, 1.045522803211208, 1.052257378424178, 1.049457613497551, 1.059166105139699, 1.044092052050142, 1.052281831771311, 1.055904741417114, 1.053115460712558, 1.062787859457203, 1.021083908580877),
    @(1.02, 1.045964505764166, 1.052618375526911, 1.049879237287132, 1.059597408428171, 1.044228383223917, 1.052592683782968, 1.056197340402387, 1.053468250195625, 1.063151031287817, 1.02119160934528),
    @(1.02, 1.047703031036519, 1.05403923742042, 1.051539511432036, 1.061295941430034, 1.044762244052666, 1.053814899590658, 1.057347532017761, 1.05485622615875, 1.064580038989138, 1.021614608036587),
    @(1.02, 1.049719699010481, 1.055687373477433, 1.053467000568463, 1.063268140253982, 1.045375848695526, 1.055229979060163, 1.058678686609494, 1.056465023323535, 1.066236775461223, 1.022103391797692)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $rowNum = $i + 2
    $rowVals = $data[$i]
    for ($j = 0; $j -lt $cols.Length; $j++) {
        $ws.Range($cols[$j] + $rowNum).Value = $rowVals[$j]
    }
}
